$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Mark individual W43 (row 6, metric "AUC, baseline") as an outlier, the
#    same way rows 4 and 7 are already flagged: highlight fill across the row
#    and put an "X" in the Remove? column (P).
# ---------------------------------------------------------------------------
$ws.Range("A4:O4").Copy()
$ws.Range("A6:O6").PasteSpecial(-4122)   # xlPasteFormats

$ws.Range("P4").Copy()
$ws.Range("P6").PasteSpecial(-4122)      # xlPasteFormats
$ws.Range("P6").Value = "X"

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 2) Remove the standalone "SlopeTime" row for individual S49 (row 19) - it is
#    no longer reported as an outlier group on its own.
# ---------------------------------------------------------------------------
$ws.Rows.Item(19).Delete()

# After the delete, the rows below have shifted up by one:
#   row19 = P22  (was "SlopeTime, DNA")
#   row20 = <blank spacer>
#   row21 = S37  (DNA)
#   row22 = W44  (DNA)
#   row23 = S47  (DNA)
#   row24 = P20  (DNA)
#
# The target layout re-sorts the "DNA" block by individual and appends P22
# (now re-labelled simply "DNA") to the end of that block:
#   row19 = <blank spacer>
#   row20 = S37  (DNA)
#   row21 = W44  (DNA)   -- de-emphasized with a grey font
#   row22 = S47  (DNA)   -- de-emphasized with a grey font
#   row23 = P20  (DNA)
#   row24 = P22  (DNA)

# Clear everything in the affected block first (values + formats) so it can
# be rebuilt cleanly in the desired order.
$ws.Range("A19:P24").Clear()

# --- row 20: DNA / S37 ------------------------------------------------------
$ws.Cells.Item(20,1).Value = "DNA"
$ws.Cells.Item(20,2).Value = "S37"
$ws.Cells.Item(20,7).Value = 56.837577269999997
$ws.Cells.Item(20,8).Value = 3982
$ws.Cells.Item(20,9).Value = "M"
$ws.Cells.Item(20,10).Value = 1
$ws.Cells.Item(20,11).Value = 3961
$ws.Cells.Item(20,12).Value = "AgStation"
$ws.Cells.Item(20,13).Value = 158
$ws.Cells.Item(20,14).Value = 35
$ws.Cells.Item(20,15).Value = "NA"

# --- row 21: DNA / W44 (grey, de-emphasized font) ---------------------------
$ws.Cells.Item(21,1).Value = "DNA"
$ws.Cells.Item(21,2).Value = "W44"
$ws.Cells.Item(21,7).Value = 181.68805750000001
$ws.Cells.Item(21,8).Value = 4173
$ws.Cells.Item(21,9).Value = "F"
$ws.Cells.Item(21,10).Value = 2
$ws.Cells.Item(21,11).Value = 4064
$ws.Cells.Item(21,12).Value = "AgStation"
$ws.Cells.Item(21,13).Value = 110.3
$ws.Cells.Item(21,14).Value = 36.5
$ws.Cells.Item(21,15).Value = 2.7
# only colour the cells that actually hold data (skip the unused C:F columns)
$ws.Range("A21:B21").Font.Color = 10921638   # grey (RGB 166,166,166)
$ws.Range("G21:O21").Font.Color = 10921638

# --- row 22: DNA / S47 (grey, de-emphasized font) ---------------------------
$ws.Cells.Item(22,1).Value = "DNA"
$ws.Cells.Item(22,2).Value = "S47"
$ws.Cells.Item(22,7).Value = 171.78411389999999
$ws.Cells.Item(22,8).Value = 4126
$ws.Cells.Item(22,9).Value = "F"
$ws.Cells.Item(22,10).Value = 2
$ws.Cells.Item(22,11).Value = 4065
$ws.Cells.Item(22,12).Value = "AgStation"
$ws.Cells.Item(22,13).Value = 162
$ws.Cells.Item(22,14).Value = 36
$ws.Cells.Item(22,15).Value = 0.6
# only colour the cells that actually hold data (skip the unused C:F columns)
$ws.Range("A22:B22").Font.Color = 10921638   # grey (RGB 166,166,166)
$ws.Range("G22:O22").Font.Color = 10921638

# --- row 23: DNA / P20 ------------------------------------------------------
$ws.Cells.Item(23,1).Value = "DNA"
$ws.Cells.Item(23,2).Value = "P20"
$ws.Cells.Item(23,7).Value = 101.1179403
$ws.Cells.Item(23,8).Value = 4130
$ws.Cells.Item(23,9).Value = "F"
$ws.Cells.Item(23,10).Value = 2
$ws.Cells.Item(23,11).Value = 4073
$ws.Cells.Item(23,12).Value = "AgStation"
$ws.Cells.Item(23,13).Value = 149.80000000000001
$ws.Cells.Item(23,14).Value = 36
$ws.Cells.Item(23,15).Value = 0.2

# --- row 24: DNA / P22 (previously "SlopeTime, DNA") -------------------------
$ws.Cells.Item(24,1).Value = "DNA"
$ws.Cells.Item(24,2).Value = "P22"
$ws.Cells.Item(24,6).Value = 3.1669999999999998
$ws.Cells.Item(24,7).Value = 10.08733518
$ws.Cells.Item(24,8).Value = 4108
$ws.Cells.Item(24,9).Value = "G"
$ws.Cells.Item(24,10).Value = 2
$ws.Cells.Item(24,11).Value = 4066
$ws.Cells.Item(24,12).Value = "AgStation"
$ws.Cells.Item(24,13).Value = 111.4
$ws.Cells.Item(24,14).Value = 39
$ws.Cells.Item(24,15).Value = 0.3

$excel.CutCopyMode = 0
